$d = $word.ActiveDocument

function Insert-XmlRuns {
    param($Range, $Texts)
    $body = ""
    foreach ($t in $Texts) {
        $escaped = $t -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
        if ($t.StartsWith(" ") -or $t.EndsWith(" ")) {
            $body += '<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
        } else {
            $body += '<w:r><w:t>' + $escaped + '</w:t></w:r>'
        }
    }
    $xmlPkg = '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $body + '</w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $null = $Range.InsertXML($xmlPkg)
}

# --- Paragraph 1: drop the stray bookmarkEnd (id 0) that sits right after
#     bookmarkStart; the matching bookmarkEnd moves to the end of the last
#     paragraph further down. Rebuild paragraph 1 keeping bookmarkStart only.
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$xmlPkg1 = '<?xml version="1.0" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:r><w:t xml:space="preserve">Oliver undertook the roles of lead QA as well as the lead web developer in this group project. As the most experienced PHP programmer in the group, the development and bulk of the implementation of the TaskerMAN component </w:t></w:r>' +
    '<w:r><w:t>fell under his jurisdiction. This includes the underlying logic that generates database queries and the presentational layer that the user navigates. He also worked closely with members in the group responsible for testing to ensure that the HTML5 and JavaScript validation was fit for purpose.</w:t></w:r>' +
    '</w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$null = $r1.InsertXML($xmlPkg1)

# --- Paragraph 2: replace the single run with seven runs of new wording.
$p2 = $d.Paragraphs.Item(2)
$texts2 = @(
    "Oliver",
    " also ensured ",
    "other members of the Web Team ",
    "were fully involved wherever possible with tasks such such as",
    " presentational work to ensure everyone had ",
    "work to complete",
    ". "
)
Insert-XmlRuns $p2.Range $texts2

# --- Insert a brand-new paragraph after paragraph 2 with the QA text.
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
$r3.Collapse(1)
$texts3 = @(
    "Oliver carried out his QA duties in meetings, working alongside the Deputy QA and Team Leader to ensure documentation and all deliverables were in accordance with the functional requirements. Their thorough, systematic approach resulted in the achievement of consistent, positive feedback. "
)
Insert-XmlRuns $r3 $texts3

# --- Final paragraph ("Punctual, hard-working..."): split the text into two
#     runs with a new ending, and append the bookmarkEnd that used to sit in
#     paragraph 1. This is the very last paragraph in the body, so its Range
#     has no trailing paragraph-mark character of its own; InsertXML-ing a
#     full <w:p> there appends an extra (phantom) empty paragraph after it.
#     Work around that by deleting the spurious mark once the real content
#     is in place, which merges the phantom paragraph back out of existence.
$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
$xmlPkg4 = '<?xml version="1.0" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' +
    '<w:r><w:t xml:space="preserve">Punctual, hard-working and sometimes a bit of a perfectionist, Oliver ensured that he attended as many necessary meetings and group work sessions as possible to ensure work was done </w:t></w:r>' +
    '<w:r><w:t>at a consistently high quality.</w:t></w:r>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '</w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$null = $r4.InsertXML($xmlPkg4)

$p4 = $d.Paragraphs.Item(4)
if ($d.Paragraphs.Count -gt 4) {
    $markRange = $d.Range($p4.Range.End - 1, $p4.Range.End)
    $null = $markRange.Delete()
}

Write-Output "Done"
